{"js": "// Word JavaScript API (Office.js) script.\n// Applies the wording tweaks from the commit \"Correccion 2, entrega 1\" to\n// the two paragraphs in the document body:\n//\n//  1) \"...tecnolog\u00edas innovadoras, que permitan...\"\n//        -> \"...tecnolog\u00edas innovadoras que permitan...\"   (drop the comma)\n//  2) \"...producci\u00f3n hidroel\u00e9ctrica, asegurando...\"\n//        -> \"...producci\u00f3n hidroel\u00e9ctrica asegurando...\"    (drop the comma)\n//  3) \"...asegurando un suministro continuo...\"\n//        -> \"...asegurando as\u00ed un suministro continuo...\"   (insert \"as\u00ed \")\n//\n// (The diff also shows a second paragraph, \"Al implementar este tipo de\n// soluciones...\", whose run is merely split around \"variables, esta\" with\n// no net text change, so there is nothing to edit there.)\n\nconst body = context.document.body;\n\nasync function replaceOnce(findText, replaceText) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + findText);\n  }\n\n  results.items[0].insertText(replaceText, \"Replace\");\n  await context.sync();\n}\n\n// 1) Drop the comma after \"innovadoras\".\nawait replaceOnce(\n  \"tecnolog\u00edas innovadoras, que permitan el almacenamiento de energ\u00eda\",\n  \"tecnolog\u00edas innovadoras que permitan el almacenamiento de energ\u00eda\"\n);\n\n// 2) Drop the comma after \"hidroel\u00e9ctrica\" (before \"asegurando\").\nawait replaceOnce(\n  \"la producci\u00f3n hidroel\u00e9ctrica, asegurando un suministro\",\n  \"la producci\u00f3n hidroel\u00e9ctrica asegurando un suministro\"\n);\n\n// 3) Insert \"as\u00ed\" after \"asegurando\".\nawait replaceOnce(\n  \"asegurando un suministro continuo de energ\u00eda\",\n  \"asegurando as\u00ed un suministro continuo de energ\u00eda\"\n);\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the wording tweaks from the commit \"Correccion 2, entrega 1\" to\n# the paragraph that starts with \"Para optimizar la hidrogeneraci\u00f3n...\":\n#\n#  1) \"...tecnolog\u00edas innovadoras, que permitan...\"\n#        -> \"...tecnolog\u00edas innovadoras que permitan...\"   (drop the comma)\n#  2) \"...producci\u00f3n hidroel\u00e9ctrica, asegurando...\"\n#        -> \"...producci\u00f3n hidroel\u00e9ctrica asegurando...\"    (drop the comma)\n#  3) \"...asegurando un suministro continuo...\"\n#        -> \"...asegurando as\u00ed un suministro continuo...\"   (insert \"as\u00ed \")\n#\n# (The diff also shows a second paragraph, \"Al implementar este tipo de\n# soluciones...\", whose run is merely split around \"variables, esta\" with\n# no net text change, so there is nothing to edit there.)\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\n# 1) Drop the comma after \"innovadoras\".\nReplace-Text \"tecnolog\u00edas innovadoras, que permitan el almacenamiento de energ\u00eda\" \"tecnolog\u00edas innovadoras que permitan el almacenamiento de energ\u00eda\"\n\n# 2) Drop the comma after \"hidroel\u00e9ctrica\" (before \"asegurando\").\nReplace-Text \"la producci\u00f3n hidroel\u00e9ctrica, asegurando un suministro\" \"la producci\u00f3n hidroel\u00e9ctrica asegurando un suministro\"\n\n# 3) Insert \"as\u00ed\" after \"asegurando\".\nReplace-Text \"asegurando un suministro continuo de energ\u00eda\" \"asegurando as\u00ed un suministro continuo de energ\u00eda\"\n"}
